$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 47: a "Lesser Hydra" bite attack entry with poison/acid special note.
$ws.Range("A47").Value = "Bite"
$ws.Range("B47").Value = 40
$ws.Range("C47").Value = "1D6"
$ws.Range("E47").Value = 4
$ws.Range("F47").Value = "Lesser Hydra"
$ws.Range("G47").Value = "Bite"
$ws.Range("H47").Value = "Poison POT is systemic and equal to POW; acid is POT equal to POW/2. Wyvern antidote is half-effective against this poison."
$ws.Range("J47").Value = "+ poison+ acid"
